# [Outlook] (preview) Add mappings for new calendar properties (Compose)
#
# Adds 8 new rows to the "Snippets" table (A1:E223 -> A1:E231) describing
# the getIsAllDayEvent / setIsAllDayEventTrue / getSensitivity /
# setSensitivityConfidential snippet mappings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$rows = @(
    @("AppointmentCompose", "isAllDayEvent", $null, "outlook-calendar-properties-apis", "getIsAllDayEvent"),
    @("IsAllDayEvent",      "getAsync",      2,     "outlook-calendar-properties-apis", "getIsAllDayEvent"),
    @("AppointmentCompose", "isAllDayEvent", $null, "outlook-calendar-properties-apis", "setIsAllDayEventTrue"),
    @("IsAllDayEvent",      "setAsync",      2,     "outlook-calendar-properties-apis", "setIsAllDayEventTrue"),
    @("AppointmentCompose", "sensitivity",   $null, "outlook-calendar-properties-apis", "getSensitivity"),
    @("Sensitivity",        "getAsync",      2,     "outlook-calendar-properties-apis", "getSensitivity"),
    @("AppointmentCompose", "sensitivity",   $null, "outlook-calendar-properties-apis", "setSensitivityConfidential"),
    @("Sensitivity",        "setAsync",      2,     "outlook-calendar-properties-apis", "setSensitivityConfidential")
)

$lastRow = 1
foreach ($row in $rows) {
    $newRow = $lo.ListRows.Add()
    $cells = $newRow.Range.Cells

    $cells.Item(1, 1).Value = $row[0]
    $cells.Item(1, 2).Value = $row[1]
    if ($row[2] -ne $null) {
        $cells.Item(1, 3).Value = $row[2]
    }
    $cells.Item(1, 4).Value = $row[3]
    $cells.Item(1, 5).Value = $row[4]

    $lastRow = $newRow.Range.Row
}

$ws.Range("A" + $lastRow).Select()
